$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QB")

$ws.Range("A4").Value = "J.Fromm"
$ws.Range("B4:L4").Value = 0

$ws.Range("L5").Select()
